# Mise à jour de l'application
# Adds the 2025-11-26 (Excel serial 45987) wellness entries for the whole
# squad to the bottom of the tracking sheet, widens column B so the
# longer player names display fully, and refreshes the selection to
# point at the newly entered rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 614
$newRows = @(
    @{ Name = "Omar Benyounes";    Volume = 70; Intensite = 5; Fatigue = 7; Douleur = 4; Loc = "Ischio";      Plaisir = 8 },
    @{ Name = "Naim Ighbane";      Volume = 70; Intensite = 3; Fatigue = 3; Douleur = 5; Loc = "Genou droit"; Plaisir = 3 },
    @{ Name = "Karim Belmahi";     Volume = 70; Intensite = 8; Fatigue = 8; Douleur = 5; Loc = "Ischio";      Plaisir = 10 },
    @{ Name = "Maé Clavel";        Volume = 70; Intensite = 5; Fatigue = 5; Douleur = 0; Loc = "";            Plaisir = 6 },
    @{ Name = "Yoann Martelat";    Volume = 70; Intensite = 5; Fatigue = 4; Douleur = 5; Loc = "Genou";       Plaisir = 6 },
    @{ Name = "Emmanuel Valey";    Volume = 70; Intensite = 7; Fatigue = 6; Douleur = 0; Loc = "";            Plaisir = 7 },
    @{ Name = "Karahali Souaré";   Volume = 70; Intensite = 6; Fatigue = 6; Douleur = 6; Loc = "Cheville";    Plaisir = 6 },
    @{ Name = "Mattheo Haon";      Volume = 70; Intensite = 6; Fatigue = 5; Douleur = 0; Loc = "";            Plaisir = 8 },
    @{ Name = "Hedi Nasri";        Volume = 70; Intensite = 7; Fatigue = 8; Douleur = 7; Loc = "Ischio";      Plaisir = 7 },
    @{ Name = "Naim Dhib";         Volume = 70; Intensite = 6; Fatigue = 7; Douleur = 4; Loc = "Cote";        Plaisir = 4 },
    @{ Name = "Ilan Ihaddadene";   Volume = 70; Intensite = 5; Fatigue = 7; Douleur = 0; Loc = "";            Plaisir = 9 }
)

# Template rows already on the sheet carry the right styles: row 607 for
# an entry with a "Localisation douleur" value, row 608 for a blank one.
$templateWithLoc = $ws.Range("A607:I607")
$templateNoLoc = $ws.Range("A608:I608")

$dateSerial = 45987
$firstNewRow = $lastRow + 1
$lastNewRow = $lastRow + $newRows.Count
$row = $lastRow
foreach ($entry in $newRows) {
    $row = $row + 1
    $destRow = $ws.Range("A" + $row + ":I" + $row)

    if ($entry.Loc -ne "") {
        $templateWithLoc.Copy()
    } else {
        $templateNoLoc.Copy()
    }
    $destRow.PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $dateSerial
    $ws.Cells.Item($row, 2).Value = $entry.Name
    $ws.Cells.Item($row, 3).Value = $entry.Volume
    $ws.Cells.Item($row, 4).Value = $entry.Intensite
    $ws.Cells.Item($row, 5).Value = $entry.Fatigue
    $ws.Cells.Item($row, 6).Value = $entry.Douleur

    if ($entry.Loc -ne "") {
        $ws.Cells.Item($row, 7).Value = $entry.Loc
    }

    $ws.Cells.Item($row, 8).Value = $entry.Plaisir
}

# Extend the existing "Charge" formula (Volume * Intensité) down through
# the new rows in one shot, same C*D pattern as the rest of the column.
$ws.Range("I" + $firstNewRow + ":I" + $lastNewRow).Formula = "=C" + $firstNewRow + "*D" + $firstNewRow

# Column B now needs to be wide enough for the longer player names.
$ws.Columns.Item(2).ColumnWidth = 14.5

# Move the selection to match the author's last on-screen position.
$ws.Range("L620").Select()
